# Updated symbol list on Sun Dec 18 15:43:05 UTC 2022 with GitHub Actions
#
# This script applies the per-row "Price" (column D) refreshes captured in the
# commit diff, plus the BKEXToken/CEJI row swap (rows 42-43, columns B/C/D/E).
#
# All "Price" values are stored as literal text in the original workbook
# (t="inlineStr"), not numbers, so every numeric-looking string is written
# with NumberFormat "@" first to keep Excel from auto-converting it to a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# --- Column D (Price) updates -------------------------------------------
Set-TextValue "D2"  "246.16"
Set-TextValue "D4"  "5.457"
Set-TextValue "D5"  "0.05614"
Set-TextValue "D6"  "6.472"
Set-TextValue "D7"  "0.8055"
Set-TextValue "D9"  "0.1432"
Set-TextValue "D10" "0.07346"
Set-TextValue "D11" "0.03199"
Set-TextValue "D12" "0.02939"
Set-TextValue "D13" "0.09262"
Set-TextValue "D14" "0.001665"
Set-TextValue "D15" "3.208"
Set-TextValue "D16" "0.04731"
Set-TextValue "D18" "0.006397"
Set-TextValue "D19" "0.001058"
Set-TextValue "D20" "0.004105"
Set-TextValue "D22" "3.979"
Set-TextValue "D23" "3.384"
Set-TextValue "D24" "2.124"
Set-TextValue "D26" "0.1278"
Set-TextValue "D27" "0.0002909"
Set-TextValue "D40" "0.04154"
Set-TextValue "D41" "0.006901"
Set-TextValue "D44" "0.009010"
Set-TextValue "D45" "0.00005661"
Set-TextValue "D47" "0.6821"
Set-TextValue "D48" "0.01748"

# --- Rows 42/43: BKEXToken and CEJI swap places -------------------------
# Row 42 becomes BKEXToken (was CEJI); row 43 becomes CEJI (was BKEXToken).
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1038"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002979"
$ws.Range("E43").Value = "42CEJICEJI"
